# Update column F (dSF) values to reflect the repulled data / mean calculation
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    2  = -3
    3  = -3
    4  = -4
    5  = 1
    7  = 4
    8  = 2
    9  = -7
    10 = 1
    11 = -5
    12 = -1
    14 = -1
    15 = 4
    17 = -2
    18 = 3
    19 = 1
    22 = -1
    23 = -1
    24 = -3
    25 = 7
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
